$wb = $excel.ActiveWorkbook

# Work on the "TourismCertificate" worksheet (maps to xl/worksheets/sheet9.xml)
$ws = $wb.Worksheets.Item("TourismCertificate")

# Insert a new blank row above the current row 3, pushing the old row 3
# (and its data) down to row 4.
$ws.Rows.Item(3).Insert()

# The newly inserted row should only keep formatting on column E (s="1"),
# the rest of the row (A3:D3) should have no cell definition at all.
$ws.Range("A3:D3").Clear()

# Make TourismCertificate the active sheet/tab and select cell A3 on it.
$ws.Activate()
$ws.Range("A3").Select() | Out-Null
